$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.042.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.960.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.77%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '327.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.46%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4973'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4206'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.89%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09235'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.03%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.096'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.981.97'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.434'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.12%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.834'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -7.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '91.49'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -8.94%  '
$ws.Range("E18").Value = '  -5.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06694'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -7.92%  '
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.937'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.068.50'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.78%  '
$ws.Range("E24").Value = '  -3.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.272'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.210.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.59'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '155.88'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.96%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.277'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.252'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -9.39%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.044'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -7.37%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09824'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.21%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.525'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -8.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.800'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -7.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.681'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.89%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02422'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -7.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.320'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.19%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.011'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.66%  '
$ws.Range("E40").Value = '  -5.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6435'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.39'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1978'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.07%  '
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6221'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -7.41%  '
$ws.Range("B46").Value = 'NEARProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.200'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -7.05%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.27'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.313'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.465'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.61%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00000000334'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.00%  '
